# "added new usb port" -- add a GPIOE pin row (row 37) to the "Pins" sheet,
# directly below the existing GPIOD row (row 36), mirroring its layout:
#   col A  -> port name   (GPIOE)
#   col B..Q -> pin names (PE0 .. PE15)
# Formatting is copied from row 36 so the new row matches the existing
# table's look (bordered cells, green "Good" fill on the pin-name cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pins")

# --- Formats first (copy from the row above so fonts/fills/borders match) ---
$ws.Range("A36").Copy()
$ws.Range("A37").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B36").Copy()
$ws.Range("B37:Q37").PasteSpecial(-4122)   # xlPasteFormats

$ws.Application.CutCopyMode = 0

# --- Values: GPIOE port header + PE0..PE15 pins ---
$ws.Range("A37").Value = "GPIOE"
$ws.Range("B37").Value = "PE0"
$ws.Range("C37").Value = "PE1"
$ws.Range("D37").Value = "PE2"
$ws.Range("E37").Value = "PE3"
$ws.Range("F37").Value = "PE4"
$ws.Range("G37").Value = "PE5"
$ws.Range("H37").Value = "PE6"
$ws.Range("I37").Value = "PE7"
$ws.Range("J37").Value = "PE8"
$ws.Range("K37").Value = "PE9"
$ws.Range("L37").Value = "PE10"
$ws.Range("M37").Value = "PE11"
$ws.Range("N37").Value = "PE12"
$ws.Range("O37").Value = "PE13"
$ws.Range("P37").Value = "PE14"
$ws.Range("Q37").Value = "PE15"

# --- Match the author's final cursor position in the saved file ---
[void]$ws.Range("F36").Select()
